# Update Document for Assets
#
# 1. Turn the plain-text law-firm website URL into a real hyperlink
#    (adds a new external relationship, e.g. rId22, and applies the
#    built-in "Hyperlink" character style to the run).
# 2. Change the "Attorney" table's preferred width from "auto" to a
#    100% ("pct", 5000 fiftieths-of-a-percent) layout width, matching
#    the Assets/Endowments table above it.

$d = $word.ActiveDocument

# --- 1) Hyperlink-ify the website text -------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("http://brunsonlawsc.com/", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $d.Hyperlinks.Add($rng, "http://brunsonlawsc.com/", "", "", "http://brunsonlawsc.com/")
}

# --- 2) Stretch the Attorney table to 100% width ----------------------
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    if ($tbl.PreferredWidthType -ne 2) {
        $tbl.PreferredWidthType = 2
        $tbl.PreferredWidth = 250
    }
}
